$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the existing "12-char" column width (as COM character units) so the
# new columns we add end up with the identical raw <col width="12".../> value
# that columns A:C already have.
$colWidth = $ws.Columns("A:A").ColumnWidth

# Strip the style (and therefore the now-unwanted custom yyyy.mm.dd number
# format) that used to be applied at the column level for A:C.
$ws.Columns("A:C").ClearFormats()

# Two distinct date/time serial values used across the row.
$dateVal1 = 41105.8449537037
$dateVal2 = 42719.461805555555

# Columns A-J, alternating date (odd) / time (even) number formats; the first
# six columns (A-F) use $dateVal1, the last four (G-J) use $dateVal2.
$cols = @("A","B","C","D","E","F","G","H","I","J")
for ($i = 0; $i -lt $cols.Count; $i++) {
    $colLetter = $cols[$i]
    $rng = $ws.Range($colLetter + "1")

    if ($i -lt 6) {
        $rng.Value = $dateVal1
    } else {
        $rng.Value = $dateVal2
    }

    if (($i % 2) -eq 0) {
        # odd columns (A, C, E, G, I) -> date format (numFmtId 14)
        $rng.NumberFormat = "mm-dd-yy"
    } else {
        # even columns (B, D, F, H, J) -> time format (numFmtId 21)
        $rng.NumberFormat = "h:mm:ss"
    }

    $ws.Columns($colLetter + ":" + $colLetter).ColumnWidth = $colWidth
}
